$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new Value entries for R1 (2k), R2 (221R), R4 (10K)
$ws.Range("D8").Value = "2k"
$ws.Range("D9").Value = "221R"
$ws.Range("D11").Value = "10K"

# Update the active sheet's selection/view
$ws.Activate()
$ws.Range("A3:I17").Select()
$excel.ActiveWindow.ActiveCell
